$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.813.17'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.369.58'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.43'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.47'
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.585'
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  -3.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.38'
$ws.Range("E13").Value = '  -3.67%  '
$ws.Range("D14").Value = '2.795.88'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '59.783.20'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '2.374.88'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.09'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '320.75'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.03'
$ws.Range("E23").Value = '  -3.74%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.40'
$ws.Range("E26").Value = '  -2.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.36'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("D29").Value = '0.0₃0758'
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.81'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.05'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  +9.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.397'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '317.99'
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '144.87'
$ws.Range("E42").Value = '  +4.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.52'
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0968'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.61'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0510'
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.569'
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.06'
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").Value = '  -1.66%  '
